$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '57.290.01'
$ws.Cells.Item(2, 5).Value = '  -0.69%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.096.11'
$ws.Cells.Item(3, 5).Value = '  -0.20%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.05%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''524.38'
$ws.Cells.Item(5, 5).Value = '  +0.07%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''136.57'
$ws.Cells.Item(6, 5).Value = '  -3.54%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''1.00'
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.092.80'
$ws.Cells.Item(8, 5).Value = '  -0.33%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +2.21%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''7.31'
$ws.Cells.Item(10, 5).Value = '  +1.31%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -1.15%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  +2.02%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '3.628.00'
$ws.Cells.Item(13, 5).Value = '  -0.26%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +2.28%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''25.20'
$ws.Cells.Item(15, 5).Value = '  -2.01%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.86%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '57.357.87'
$ws.Cells.Item(17, 5).Value = '  -0.73%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.089.16'
$ws.Cells.Item(18, 5).Value = '  -0.41%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -2.63%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''12.33'
$ws.Cells.Item(20, 5).Value = '  -3.69%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  -2.57%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''344.66'
$ws.Cells.Item(22, 5).Value = '  +1.33%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.02%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''67.55'
$ws.Cells.Item(24, 5).Value = '  +1.02%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -2.70%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''0.166'
$ws.Cells.Item(26, 5).Value = '  -2.28%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''1.00'
$ws.Cells.Item(27, 5).Value = '  -0.02%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -3.50%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -0.07%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''7.33'
$ws.Cells.Item(30, 5).Value = '  +2.10%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''1.87'
$ws.Cells.Item(31, 5).Value = '  -0.20%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''6.00'
$ws.Cells.Item(32, 5).Value = '  -7.63%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''20.70'
$ws.Cells.Item(33, 5).Value = '  -1.17%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''4.90'
$ws.Cells.Item(34, 5).Value = '  +6.08%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -4.31%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''158.97'
$ws.Cells.Item(36, 5).Value = '  +1.95%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''6.04'
$ws.Cells.Item(37, 5).Value = '  -1.50%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''25.72'
$ws.Cells.Item(38, 5).Value = '  -5.51%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -1.97%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).Value = '''1.60'
$ws.Cells.Item(40, 5).Value = '  +5.56%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Hedera'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(41, 4).Value = '''0.0656'
$ws.Cells.Item(41, 5).Value = '  -0.95%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''4.08'
$ws.Cells.Item(42, 5).Value = '  +3.29%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.696'
$ws.Cells.Item(43, 5).Value = '  +1.99%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Maker'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(44, 4).Value = '2.369.26'
$ws.Cells.Item(44, 5).Value = '  +3.23%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'OKB'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(45, 4).Value = '''36.48'
$ws.Cells.Item(45, 5).Value = '  -0.81%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(46, 4).Value = '''0.999'
$ws.Cells.Item(46, 5).Value = '  -0.06%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(47, 4).Value = '''0.0265'
$ws.Cells.Item(47, 5).Value = '  +2.15%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'ONDO'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(48, 4).Value = '''0.967'
$ws.Cells.Item(48, 5).Value = '  -1.49%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Cosmos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(49, 4).Value = '''5.94'
$ws.Cells.Item(49, 5).Value = '  -1.37%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(50, 4).Value = '''19.61'
$ws.Cells.Item(50, 5).Value = '  -4.57%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'SuiNetwork'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(51, 4).Value = '''0.753'
$ws.Cells.Item(51, 5).Value = '  +2.60%  '
